$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2204.2727
$ws.Range("I38").Value = 1748.2858
$ws.Range("J38").Value = 3002.25
$ws.Range("K38").Value = 5244.857400000001
$ws.Range("L38").Value = 9006.75
$ws.Range("M38").Value = -4872.857400000001
$ws.Range("N38").Value = -9750.75
$ws.Range("H39").Value = 926.2
$ws.Range("I39").Value = 925.61536
$ws.Range("J39").Value = 930
$ws.Range("K39").Value = 2776.84608
$ws.Range("L39").Value = 2790
$ws.Range("M39").Value = -2480.84608
$ws.Range("N39").Value = -3382
$ws.Range("H40").Value = 2626.1333
$ws.Range("I40").Value = 2561.375
$ws.Range("J40").Value = 2700.1428
$ws.Range("K40").Value = 2561.375
$ws.Range("L40").Value = 2700.1428
$ws.Range("M40").Value = -2386.375
$ws.Range("N40").Value = -3050.1428
$ws.Range("H43").Value = 612.1667
$ws.Range("I43").Value = 1100.5
$ws.Range("J43").Value = 368
$ws.Range("K43").Value = 1100.5
$ws.Range("L43").Value = 368
$ws.Range("M43").Value = -1031.5
$ws.Range("N43").Value = -506
$ws.Range("H100").Value = 15152861
$ws.Range("I100").Value = 27779378
$ws.Range("J100").Value = 1040
$ws.Range("K100").Value = 27779378
$ws.Range("L100").Value = 1040
$ws.Range("M100").Value = -27778837
$ws.Range("N100").Value = -2122
$ws.Range("H111").Value = 1181.4166
$ws.Range("I111").Value = 1078
$ws.Range("K111").Value = 3234
$ws.Range("M111").Value = -167
$ws.Range("H132").Value = 303105.47
$ws.Range("I132").Value = 347897.84
$ws.Range("K132").Value = 1043693.52
$ws.Range("M132").Value = -1041163.52
$ws.Range("H138").Value = 4363226
$ws.Range("I138").Value = 1264826.5
$ws.Range("J138").Value = 6292418.5
$ws.Range("K138").Value = 3794479.5
$ws.Range("L138").Value = 18877255.5
$ws.Range("M138").Value = -3789339.5
$ws.Range("N138").Value = -18887535.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 14971.5
$ws.Range("J112").Value = 14971.5
$ws.Range("L112").Value = 14971.5
$ws.Range("N112").Value = -17925.5
$ws.Range("H114").Value = 21899.334
$ws.Range("J114").Value = 21899.334
$ws.Range("L114").Value = 21899.334
$ws.Range("N114").Value = -30577.334
$ws.Range("H124").Value = 33681
$ws.Range("J124").Value = 33681
$ws.Range("L124").Value = 33681
$ws.Range("N124").Value = -43501
$ws.Range("H132").Value = 2355.6296
$ws.Range("I132").Value = 1781.4736
$ws.Range("J132").Value = 3719.25
$ws.Range("K132").Value = 5344.4208
$ws.Range("L132").Value = 11157.75
$ws.Range("M132").Value = -2814.4208
$ws.Range("N132").Value = -16217.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H99").Value = 1787.5
$ws.Range("I99").Value = 1400
$ws.Range("J99").Value = 2115.3845
$ws.Range("K99").Value = 1400
$ws.Range("L99").Value = 2115.3845
$ws.Range("M99").Value = 98
$ws.Range("N99").Value = -5111.3845
$ws.Range("H107").Value = 794.5714
$ws.Range("I107").Value = 794.5714
$ws.Range("K107").Value = 794.5714
$ws.Range("M107").Value = 1125.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1567.85
$ws.Range("I31").Value = 991.0714
$ws.Range("K31").Value = 991.0714
$ws.Range("M31").Value = -696.0714
$ws.Range("H34").Value = 1567.85
$ws.Range("I34").Value = 991.0714
$ws.Range("K34").Value = 991.0714
$ws.Range("M34").Value = -789.0714
$ws.Range("H58").Value = 1947
$ws.Range("I58").Value = 655.875
$ws.Range("J58").Value = 3162.1765
$ws.Range("K58").Value = 655.875
$ws.Range("L58").Value = 3162.1765
$ws.Range("M58").Value = -452.875
$ws.Range("N58").Value = -3568.1765
$ws.Range("H62").Value = 33613.715
$ws.Range("J62").Value = 6949
$ws.Range("L62").Value = 6949
$ws.Range("N62").Value = -8197
$ws.Range("H65").Value = 33613.715
$ws.Range("J65").Value = 6949
$ws.Range("L65").Value = 34745
$ws.Range("N65").Value = -40985
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H136").Value = 1947
$ws.Range("I136").Value = 655.875
$ws.Range("J136").Value = 3162.1765
$ws.Range("K136").Value = 1967.625
$ws.Range("L136").Value = 9486.529500000001
$ws.Range("M136").Value = 582.375
$ws.Range("N136").Value = -14586.5295

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 10870353
$ws.Range("I113").Value = 580.8
$ws.Range("J113").Value = 19231716
$ws.Range("K113").Value = 1742.4
$ws.Range("L113").Value = 57695148
$ws.Range("M113").Value = 427.6000000000001
$ws.Range("N113").Value = -57699488
$ws.Range("H129").Value = 2868
$ws.Range("I129").Value = 2616
$ws.Range("J129").Value = 3120
$ws.Range("K129").Value = 7848
$ws.Range("L129").Value = 9360
$ws.Range("M129").Value = -2848
$ws.Range("N129").Value = -19360
$ws.Range("H132").Value = 8334847
$ws.Range("I132").Value = 445
$ws.Range("J132").Value = 10418448
$ws.Range("K132").Value = 4005
$ws.Range("L132").Value = 93766032
$ws.Range("M132").Value = -1475
$ws.Range("N132").Value = -93771092
$ws.Range("H137").Value = 5320197.5
$ws.Range("I137").Value = 11114528
$ws.Range("J137").Value = 105299.9
$ws.Range("K137").Value = 33343584
$ws.Range("L137").Value = 315899.7
$ws.Range("M137").Value = -33338484
$ws.Range("N137").Value = -326099.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1918.75
$ws.Range("I46").Value = 1130
$ws.Range("J46").Value = 3233.3333
$ws.Range("K46").Value = 1130
$ws.Range("L46").Value = 3233.3333
$ws.Range("M46").Value = -942
$ws.Range("N46").Value = -3609.3333
$ws.Range("H61").Value = 8531.272000000001
$ws.Range("I61").Value = 8943
$ws.Range("J61").Value = 7433.3335
$ws.Range("K61").Value = 8943
$ws.Range("L61").Value = 7433.3335
$ws.Range("M61").Value = -8741
$ws.Range("N61").Value = -7837.3335
$ws.Range("H110").Value = 20001
$ws.Range("J110").Value = 20001
$ws.Range("L110").Value = 20001
$ws.Range("N110").Value = -28181
$ws.Range("H113").Value = 8531.272000000001
$ws.Range("I113").Value = 8943
$ws.Range("J113").Value = 7433.3335
$ws.Range("K113").Value = 8943
$ws.Range("L113").Value = 7433.3335
$ws.Range("M113").Value = -6773
$ws.Range("N113").Value = -11773.3335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 10229.6
$ws.Range("I45").Value = 4570
$ws.Range("J45").Value = 11644.5
$ws.Range("K45").Value = 4570
$ws.Range("L45").Value = 11644.5
$ws.Range("M45").Value = -4079
$ws.Range("N45").Value = -12626.5
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H96").Value = 83334140
$ws.Range("I96").Value = 142857760
$ws.Range("J96").Value = 1081.6
$ws.Range("K96").Value = 142857760
$ws.Range("L96").Value = 1081.6
$ws.Range("M96").Value = -142856387
$ws.Range("N96").Value = -3827.6
